$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value2 = "Discover 1 employee"
$ws.Range("J2").Value2 = 6704018

# Row 6
$ws.Range("D6").Value2 = "View all 659 employees"
$ws.Range("J6").Value2 = 2273792

# Row 9
$ws.Range("D9").Value2 = "Access all 3,353 employees"
$ws.Range("J9").Value2 = 2444446

# Row 14
$ws.Range("D14").Value2 = "Access all 120 employees"
$ws.Range("J14").Value2 = 7615115

# Row 15
$ws.Range("D15").Value2 = "Discover all 3,795 employees"
$ws.Range("J15").Value2 = 1856460

# Row 16
$ws.Range("D16").Value2 = "View all 2,043 employees"
$ws.Range("J16").Value2 = 5978459

# Row 21
$ws.Range("D21").Value2 = "Discover all 10 employees"
$ws.Range("J21").Value2 = 6230607

# Row 22
$ws.Range("D22").Value2 = "Discover all 1,723 employees"
$ws.Range("J22").Value2 = 9770839

# Row 24
$ws.Range("D24").Value2 = "View all 20 employees"
$ws.Range("J24").Value2 = 9766198

# Row 25
$ws.Range("D25").Value2 = "Discover all 27,036 employees"
$ws.Range("J25").Value2 = 1255671

# Row 28
$ws.Range("D28").Value2 = "Discover all 3,145 employees"
$ws.Range("J28").Value2 = 9277902

# Row 29
$ws.Range("J29").Value2 = 7681016

# Row 30
$ws.Range("J30").Value2 = 7450175

# Row 31
$ws.Range("D31").Value2 = "Access all 2,378 employees"
$ws.Range("J31").Value2 = 7225893

# Row 32
$ws.Range("D32").Value2 = "Access all 1,910 employees"
$ws.Range("J32").Value2 = 5475998

# Row 34
$ws.Range("D34").Value2 = "Discover all 722 employees"
$ws.Range("J34").Value2 = 9700193

# Row 35
$ws.Range("D35").Value2 = "View all 37 employees"
$ws.Range("J35").Value2 = 6134151

# Row 36
$ws.Range("B36").Value2 = "N/A"
$ws.Range("C36").Value2 = "N/A"
$ws.Range("D36").Value2 = "N/A"
$ws.Range("E36").Value2 = "N/A"
$ws.Range("G36").Value2 = "N/A"
$ws.Range("H36").Value2 = "N/A"
$ws.Range("J36").Value2 = "N/A"
$ws.Range("K36").Value2 = "N/A"

# Row 37
$ws.Range("D37").Value2 = "View all 462 employees"
$ws.Range("J37").Value2 = 2812787

# Row 38
$ws.Range("D38").Value2 = "Discover all 1,885 employees"
$ws.Range("J38").Value2 = 5943450

# Row 40
$ws.Range("B40").Value2 = "Houston, Texas"
$ws.Range("C40").Value2 = "Oil and Gas"
$ws.Range("D40").Value2 = "View all 3,171 employees"
$ws.Range("E40").Value2 = "1,001-5,000 employees"
$ws.Range("G40").Value2 = "https://www.mrcglobal.com"
$ws.Range("H40").Value2 = "mrcglobal.com"
$ws.Range("J40").Value2 = 1875436
$ws.Range("K40").Value2 = "Suite 2300"

# Row 43
$ws.Range("J43").Value2 = 4342130

# Row 45
$ws.Range("B45").Value2 = "Houston, TX"
$ws.Range("C45").Value2 = "Packaging and Containers Manufacturing"
$ws.Range("D45").Value2 = "Discover all 1,324 employees"
$ws.Range("E45").Value2 = "1,001-5,000 employees"
$ws.Range("G45").Value2 = "http://www.victorypackaging.com"
$ws.Range("H45").Value2 = "victorypackaging.com"
$ws.Range("J45").Value2 = 3465758
$ws.Range("K45").Value2 = "Ste. 1400"

# Row 47
$ws.Range("D47").Value2 = "Discover all 411 employees"
$ws.Range("J47").Value2 = 8669392

# Row 48
$ws.Range("D48").Value2 = "Access all 406 employees"
$ws.Range("J48").Value2 = 3709035

# Row 50
$ws.Range("D50").Value2 = "Discover all 445 employees"
$ws.Range("J50").Value2 = 5081958

# Row 51
$ws.Range("D51").Value2 = "View all 1,899 employees"
$ws.Range("J51").Value2 = 1283805

# Row 52
$ws.Range("D52").Value2 = "Access all 162 employees"
$ws.Range("J52").Value2 = 3823190

# Row 53
$ws.Range("D53").Value2 = "Access all 69 employees"
$ws.Range("J53").Value2 = 7977563

# Row 54
$ws.Range("D54").Value2 = "View all 603 employees"
$ws.Range("J54").Value2 = 2617057

# Row 55
$ws.Range("D55").Value2 = "View all 2,171 employees"
$ws.Range("J55").Value2 = 2097957

# Row 58
$ws.Range("D58").Value2 = "View all 982 employees"
$ws.Range("J58").Value2 = 1400464

# Row 60
$ws.Range("D60").Value2 = "View all 5,924 employees"
$ws.Range("J60").Value2 = 3464691

# Row 63
$ws.Range("B63").Value2 = "N/A"
$ws.Range("C63").Value2 = "N/A"
$ws.Range("D63").Value2 = "N/A"
$ws.Range("E63").Value2 = "N/A"
$ws.Range("G63").Value2 = "N/A"
$ws.Range("H63").Value2 = "N/A"
$ws.Range("J63").Value2 = "N/A"
$ws.Range("K63").Value2 = "N/A"

# Row 65
$ws.Range("J65").Value2 = 2027737

# Row 66
$ws.Range("J66").Value2 = 9441682

# Row 70
$ws.Range("B70").Value2 = "Dallas, TX"
$ws.Range("C70").Value2 = "Wholesale Building Materials"
$ws.Range("D70").Value2 = "Discover all 1,836 employees"
$ws.Range("E70").Value2 = "5,001-10,000 employees"
$ws.Range("G70").Value2 = "http://www.obe.com/"
$ws.Range("H70").Value2 = "obe.com/"
$ws.Range("J70").Value2 = 3341381
$ws.Range("K70").Value2 = "Suite 1050"

# Row 71
$ws.Range("D71").Value2 = "View all 3,695 employees"
$ws.Range("J71").Value2 = 8050882

# Column width adjustments (D: 30 -> 31, G: 37 -> 36)
$ws.Columns.Item(4).ColumnWidth = 31 - 5/6
$ws.Columns.Item(7).ColumnWidth = 36 - 5/6
